$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: change from the "Kattfotslav" record to the "Granbarkgnagare" record
$ws.Range("A4").Value = 130645210
$ws.Range("B4").Value = 4779
$ws.Range("E4").Value = 102306
$ws.Range("F4").Value = "Granbarkgnagare"
$ws.Range("G4").Value = "Microbregma emarginatum"
$ws.Range("H4").Value = "(Duftschmid, 1825)"
$ws.Range("Q4").Value = 667485
$ws.Range("R4").Value = 6693332
$ws.Range("AC4").ClearContents()

# Row 5: change from the "Granbarkgnagare" record to the "Kattfotslav" record
$ws.Range("A5").Value = 130645206
$ws.Range("B5").Value = 75349
$ws.Range("E5").Value = 6426
$ws.Range("F5").Value = "Kattfotslav"
$ws.Range("G5").Value = "Felipes leucopellaeus"
$ws.Range("H5").Value = "(Ach.) Frisch & G.Thor"
$ws.Range("Q5").Value = 667472
$ws.Range("R5").Value = 6693360
$ws.Range("AC5").Value = "Senvuxen"
